# Adjust property of scene: update CamOffestPos / CamOffestRot values for the
# "villageScene" row (row 2) and the "City" row (row 6), then move the
# selection/view the way the author left the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (villageScene): CamOffestPos / CamOffestRot
$ws.Range("J2").Value = "0,8,7"
$ws.Range("K2").Value = "45,180"

# Row 6 (City): CamOffestPos / CamOffestRot
$ws.Range("J6").Value = "0,8,-7"
$ws.Range("K6").Value = "45,0"

# Scroll the view so column E is the left-most visible column, and leave the
# active cell/selection on K7.
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("K7").Select()
